$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.967.49"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "  -2.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.577.84"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = "  -4.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.88"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = "  -1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.09"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = "  -2.05%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +2.43%  "

$ws.Range("E9").Value = "  -0.78%  "

$ws.Range("E10").Value = "  -1.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.57"
$ws.Range("D11").ClearFormats()

$ws.Range("E11").Value = "  +2.92%  "

$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.034.15"
$ws.Range("D13").ClearFormats()

$ws.Range("E13").Value = "  -4.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.66"
$ws.Range("D14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "61.879.09"
$ws.Range("D15").ClearFormats()

$ws.Range("E15").Value = "  -2.26%  "

$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.580.93"
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = "  -4.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.63"
$ws.Range("D18").ClearFormats()

$ws.Range("E18").Value = "  -3.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.56"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = "  -0.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.54"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = "  -1.58%  "

$ws.Range("E21").Value = "  -4.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").Value = "  +0.45%  "

$ws.Range("E23").Value = "  -2.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.52"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = "  -0.49%  "

$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.19"
$ws.Range("D27").ClearFormats()

$ws.Range("E27").Value = "  +0.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.35"
$ws.Range("D28").ClearFormats()

$ws.Range("E28").Value = "  +4.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0840"
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = "  -2.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.34"
$ws.Range("D30").ClearFormats()

$ws.Range("E31").Value = "  -2.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.80"
$ws.Range("D32").ClearFormats()

$ws.Range("E32").Value = "  -1.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.87"
$ws.Range("D33").ClearFormats()

$ws.Range("E33").Value = "  +1.68%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.43"
$ws.Range("D35").ClearFormats()

$ws.Range("E35").Value = "  +0.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.24"
$ws.Range("D36").ClearFormats()

$ws.Range("E36").Value = "  -1.62%  "

$ws.Range("E37").Value = "  +1.15%  "

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "RenderToken"
$ws.Range("B38").ClearFormats()

$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C38").ClearFormats()

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.08"
$ws.Range("D38").ClearFormats()

$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Bittensor"
$ws.Range("B39").ClearFormats()

$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C39").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "329.24"
$ws.Range("D39").ClearFormats()

$ws.Range("E39").Value = "  -3.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.915"
$ws.Range("D40").ClearFormats()

$ws.Range("E40").Value = "  -3.83%  "

$ws.Range("E41").Value = "  +0.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.59"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = "  -1.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.03"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = "  +1.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("D44").ClearFormats()

$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.607"
$ws.Range("D45").ClearFormats()

$ws.Range("E45").Value = "  -2.48%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Maker"
$ws.Range("B46").ClearFormats()

$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C46").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.118.65"
$ws.Range("D46").ClearFormats()

$ws.Range("E46").Value = "  +0.78%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Hedera"
$ws.Range("B47").ClearFormats()

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C47").ClearFormats()

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0549"
$ws.Range("D47").ClearFormats()

$ws.Range("E47").Value = "  -2.66%  "

$ws.Range("E48").Value = "  -1.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.57"
$ws.Range("D49").ClearFormats()

$ws.Range("E49").Value = "  -3.81%  "

$ws.Range("E50").Value = "  -0.54%  "

$ws.Range("E51").Value = "  -0.92%  "
